# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $val
}

$ws.Range("D2").Value = "37.307.94"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "2.050.01"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue "D5" "229.76"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("E6").Value = "  -1.59%  "
$ws.Range("E7").Value = "  +0.02%  "
Set-TextValue "D8" "56.92"
$ws.Range("E8").Value = "  -2.77%  "
$ws.Range("E9").Value = "  -1.55%  "
Set-TextValue "D10" "0.0786"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("E11").Value = "  -2.02%  "
Set-TextValue "D12" "14.92"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").Value = "2.355.19"
$ws.Range("E13").Value = "  -1.13%  "
Set-TextValue "D14" "20.85"
$ws.Range("E14").Value = "  -1.33%  "
Set-TextValue "D15" "0.757"
$ws.Range("E15").Value = "  -2.62%  "
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").Value = "2.059.87"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "37.202.90"
$ws.Range("E18").Value = "  -1.39%  "
Set-TextValue "D19" "6.09"
$ws.Range("E19").Value = "  -0.75%  "
Set-TextValue "D20" "69.50"
$ws.Range("D21").Value = "0.0₃0827"
$ws.Range("E21").Value = "  -1.73%  "
Set-TextValue "D22" "227.03"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("E23").Value = "  -0.04%  "
Set-TextValue "D24" "2.39"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  -4.93%  "
Set-TextValue "D26" "9.75"
$ws.Range("E26").Value = "  +0.25%  "
Set-TextValue "D27" "166.37"
$ws.Range("E27").Value = "  -3.16%  "
$ws.Range("E28").Value = "  -6.80%  "
Set-TextValue "D29" "19.06"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("E30").Value = "  -3.42%  "
$ws.Range("E32").Value = "  -3.77%  "
Set-TextValue "D33" "0.0619"
$ws.Range("E33").Value = "  -2.32%  "
$ws.Range("E34").Value = "  -1.93%  "
Set-TextValue "D35" "2.48"
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  -4.13%  "
Set-TextValue "D39" "5.25"
$ws.Range("E39").Value = "  -3.54%  "
$ws.Range("E40").Value = "  -5.61%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D41" "17.16"
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.494.22"
$ws.Range("E42").Value = "  +3.23%  "
Set-TextValue "D43" "2.90"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("E44").Value = "  -3.17%  "
Set-TextValue "D45" "96.76"
$ws.Range("E45").Value = "  -3.91%  "
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("E47").Value = "  -3.66%  "
Set-TextValue "D48" "3.95"
$ws.Range("E48").Value = "  -3.96%  "
Set-TextValue "D49" "7.13"
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("D51").Value = "2.239.81"
$ws.Range("E51").Value = "  -1.23%  "
